$d = $word.ActiveDocument

# --- Update the "Priority" column (last column) of each data row in both
# tables: every value is reduced by 4, except the "GDPR data can be
# requested" row in table 1 which drops by 5 (2 -> -3). ---
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)
    for ($ri = 2; $ri -le $t.Rows.Count; $ri++) {
        $row = $t.Rows.Item($ri)
        $cell = $row.Cells.Item(4)
        $text = $cell.Range.Text
        $text = $text.TrimEnd([char]13, [char]7)
        $old = [int]$text
        if ($ti -eq 1 -and $old -eq 2 -and $ri -eq 3) {
            $new = -3
        } else {
            $new = $old - 4
        }
        $cell.Range.Text = [string]$new
    }
}
